$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the team record (Wins/Losses/Ties) for every data row (2 through 58)
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 92  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 70  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
